$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 101, pushing old rows 101-108 down to 102-109.
$ws.Rows.Item(101).Insert()

# Populate the new row 101 with the new record (copy of the old row 101's
# data, but with an updated date, quality ("Calidad") and volume).
$ws.Cells.Item(101, 1).Value = 10
$ws.Cells.Item(101, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(101, 3).Value = "La Araucanía"
$ws.Cells.Item(101, 4).Value = 45127
$ws.Cells.Item(101, 5).Value = 9
$ws.Cells.Item(101, 6).Value = "Fruta"
$ws.Cells.Item(101, 7).Value = 100108
$ws.Cells.Item(101, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(101, 9).Value = 100108004
$ws.Cells.Item(101, 10).Value = "Papaya"
$ws.Cells.Item(101, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(101, 12).Value = "Primera"
$ws.Cells.Item(101, 13).Value = 45
$ws.Cells.Item(101, 14).Value = 25000
$ws.Cells.Item(101, 15).Value = 25000
$ws.Cells.Item(101, 16).Value = 25000
$ws.Cells.Item(101, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(101, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(101, 19).Value = 2500
$ws.Cells.Item(101, 20).Value = 10
